# "Add files via upload" -- appends two new transaction rows (20 & 21) to
# SheetName1, highlights two existing cells (T9, S10) in bold-on-yellow, and
# leaves the selection on the newly added rows, matching the uploaded file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Highlight two already-existing cells: bold font + yellow fill.
#    (S10 keeps its date/time number format; T9 is a plain number.)
# ---------------------------------------------------------------------
$ws.Range("S10").Font.Bold = $true
$ws.Range("S10").Interior.Color = 65535

$ws.Range("T9").Font.Bold = $true
$ws.Range("T9").Interior.Color = 65535

# ---------------------------------------------------------------------
# Helper: write a value into a cell as TEXT (shared string) even when it
# looks like a pure number (e.g. "194224"), without leaving any lingering
# cell formatting behind. A quoted formula always evaluates to a string;
# copy/paste-values then drops the formula and leaves a plain text
# constant behind, with no number-format side effects on the cell style.
# ---------------------------------------------------------------------
function Set-TextValue($rng, [string]$text) {
    $rng.Formula = "=""" + $text + """"
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------
# 2) New row 20 -- duplicate of the "Jayden Cooper" deposit (row 2), with
#    its own FTD/USD Amount and AssignedDate values.
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "AC2036750"
$ws.Range("B20").Value = "Richard"
$ws.Range("C20").Value = "Miller"
$ws.Range("D20").Value = "Jayden Cooper"
$ws.Range("E20").Value = "Deposit"
$ws.Range("F20").Value = 805.14
$ws.Range("G20").Value = "S ENG RET T1"
$ws.Range("H20").Value = "Approved"
$ws.Range("I20").Value = $false
$ws.Range("J20").Value = "P1641IN-35"
$ws.Range("K20").Value = 45443.577546296299
$ws.Range("K20").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("L20").Value = "Real"
$ws.Range("M20").Value = $false
$ws.Range("N20").Value = "Crypto"
$ws.Range("O20").Value = "Jayden Cooper"
$ws.Range("P20").Value = "ETH"
$ws.Range("R20").Value = "0x8ED71f780dF450273051c568612e476406A5C9E7"
$ws.Range("S20").Value = 45432.725555555553
$ws.Range("S20").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("T20").Value = 875.29179999999997
Set-TextValue $ws.Range("U20") "194224"
$ws.Range("V20").Value = 973
$ws.Range("W20").Value = "United Kingdom of Great Britain and Northern Ireland"
$ws.Range("X20").Value = "EUR"
$ws.Range("Y20").Value = 805.14
$ws.Range("Z20").Value = "FTD"
$ws.Range("AB20").Value = 45439.598854166667
$ws.Range("AB20").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# ---------------------------------------------------------------------
# 3) New row 21 -- new account (Caryn Buys) under the "George Angelo" FTD.
# ---------------------------------------------------------------------
$ws.Range("A21").Value = "AC2028053"
$ws.Range("B21").Value = "Caryn"
$ws.Range("C21").Value = "Buys"
$ws.Range("D21").Value = "George Angelo"
$ws.Range("E21").Value = "Deposit"
$ws.Range("F21").Value = 2590.6
$ws.Range("G21").Value = "S ENG RET T1"
$ws.Range("H21").Value = "Approved"
$ws.Range("I21").Value = $false
$ws.Range("J21").Value = "SG32IN-42"
$ws.Range("K21").Value = 45443.603715277779
$ws.Range("K21").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("L21").Value = "Real"
$ws.Range("M21").Value = $false
$ws.Range("N21").Value = "Crypto"
$ws.Range("O21").Value = "George Angelo"
$ws.Range("P21").Value = "ETH"
$ws.Range("R21").Value = "0x4576a940621806424AcCACc11fB2f3Be1e745981"
$ws.Range("S21").Value = 45419.434166666666
$ws.Range("S21").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("T21").Value = 2590.6
Set-TextValue $ws.Range("U21") "194230"
$ws.Range("V21").Value = 973
$ws.Range("W21").Value = "South Africa"
$ws.Range("X21").Value = "USD"
$ws.Range("Y21").Value = 8998.8468191299999
$ws.Range("Z21").Value = "FTD"
$ws.Range("AB21").Value = 45439.599641203706
$ws.Range("AB21").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# ---------------------------------------------------------------------
# 4) Leave the selection the way the uploaded file has it.
# ---------------------------------------------------------------------
[void]$ws.Range("V19:V21").Select()
